# Daily Status Tracker update
# - Unhide all rows on the Tracker sheet and clear the column-D (Owner)
#   AutoFilter criteria (keep the AutoFilter arrows, just show everything).
# - Reassign two task owners (row 9 and row 14) from Debashish -> Sayan / Rahul.
# - Add a new "Notes" worksheet at the end of the workbook containing the
#   Anodiam website task breakdown notes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Tracker sheet: clear the existing AutoFilter criteria and unhide rows
# ---------------------------------------------------------------------------
$tracker = $wb.Worksheets.Item("Tracker")

# Drop the current filter (column D == "Sayan") and show every row again.
$tracker.AutoFilterMode = $false
$tracker.Range("A1:H32").AutoFilter()

# Rows were hidden by the old filter - make them all visible again.
$tracker.Range("A1:A32").EntireRow.Hidden = $false

# ---------------------------------------------------------------------------
# 2. Tracker sheet: update task owners
# ---------------------------------------------------------------------------
$tracker.Range("D9").Value = "Sayan"
$tracker.Range("D14").Value = "Rahul"

# ---------------------------------------------------------------------------
# 3. Add the new "Notes" worksheet at the end with the Anodiam website notes
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$notesSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$notesSheet.Name = "Notes"

$notes = @(
    "1. Design the corresponding components for https://www.google.business.anodiam.site",
    "    A. Images & Logos",
    "    B. Videos, Youtube content. Create an youtube channel and start putting in our contents",
    "    C. Fonts (Oxygen) and font sizes",
    "    D. Color schemes (RGBA values) of each component",
    "    E. Hyperlinks and Button Click Events including whatsapp integration",
    "    F. Texts and headings (verbose properly / a lot is already there in ppts already)",
    "    G. About Us & Service should have following sections:",
    "       i.     AI Junior & Robotics",
    "      ii.     AI Professional & IT",
    "     iii.     Other Swim lanes (You may want to avoid this for now)",
    "2. Antech to develop https://www.google.business.anodiam.site using above components",
    "3. Ensure we have a button to:",
    "    A. download anodiam app from google playstore if anodiam app is not already installed",
    "    B. Use onclick event of same button to open anodiam app if it is already installed",
    "4. Test the responsiveness of https://www.google.business.sample.site",
    "5. Must have a button `"Visit anodiam.com website`": upon clicking it will open anodiam.com >> classplus/YOWWA website"
)

for ($i = 0; $i -lt $notes.Length; $i++) {
    $row = $i + 2
    $notesSheet.Cells.Item($row, 3).Value = $notes[$i]
}
